$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update INDICATOR_ID (column A) from 4 to 57 for all data rows (2-190)
$ws.Range("A2:A190").Value = 57

# Row heights recalculated (wrapped-text comment column reflow) for a few rows
$ws.Rows.Item(120).RowHeight = 75
$ws.Rows.Item(168).RowHeight = 75
$ws.Rows.Item(174).RowHeight = 60

# Restore the default view (no frozen/scrolled top-left cell) and zoom to 100%
$excel.ActiveWindow.Zoom = 100

# Update the view selection to E3 (matches the committed worksheet state)
$ws.Range("E3").Select() | Out-Null
